$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "SCRIPT/D29P11A/enter03.ssb",
    "SCRIPT/D31P11A/enter03.ssb",
    "SCRIPT/D73P11A/enter02.ssb",
    "SCRIPT/D73P22A/enter02.ssb",
    "SCRIPT/D73P24A/enter03.ssb",
    "SCRIPT/D73P26A/enter03.ssb",
    "SCRIPT/D73P28A/enter02.ssb"
)

$startRow = 14
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

$ws.Range("C6").Select()
